# NIT-9009772770.xlsx — "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The previous 37 monthly "Estado de Cuenta" periods (1703..2003, one per
# worker row 16..52) are replaced by a new set of periods that runs the
# other direction (2003 down to 1703), and the "Valor Mora" amount for each
# row flips between the two values used in this sheet (31249 / 29509) to
# match the new period each row now represents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New period label (col E) + new "Valor Mora" (col F) for every data row.
$rows = @(
    @{Row=16; Period="2003"; Mora=31249},
    @{Row=17; Period="2002"; Mora=31249},
    @{Row=18; Period="2001"; Mora=31249},
    @{Row=19; Period="1912"; Mora=31249},
    @{Row=20; Period="1911"; Mora=31249},
    @{Row=21; Period="1910"; Mora=31249},
    @{Row=22; Period="1909"; Mora=31249},
    @{Row=23; Period="1908"; Mora=31249},
    @{Row=24; Period="1907"; Mora=31249},
    @{Row=25; Period="1906"; Mora=31249},
    @{Row=26; Period="1905"; Mora=31249},
    @{Row=27; Period="1904"; Mora=31249},
    @{Row=28; Period="1903"; Mora=31249},
    @{Row=29; Period="1902"; Mora=31249},
    @{Row=30; Period="1901"; Mora=31249},
    @{Row=31; Period="1812"; Mora=31249},
    @{Row=32; Period="1811"; Mora=31249},
    @{Row=33; Period="1810"; Mora=31249},
    @{Row=34; Period="1809"; Mora=31249},
    @{Row=35; Period="1808"; Mora=29509},
    @{Row=36; Period="1807"; Mora=29509},
    @{Row=37; Period="1806"; Mora=29509},
    @{Row=38; Period="1805"; Mora=29509},
    @{Row=39; Period="1804"; Mora=29509},
    @{Row=40; Period="1803"; Mora=29509},
    @{Row=41; Period="1802"; Mora=29509},
    @{Row=42; Period="1801"; Mora=29509},
    @{Row=43; Period="1712"; Mora=29509},
    @{Row=44; Period="1711"; Mora=29509},
    @{Row=45; Period="1710"; Mora=29509},
    @{Row=46; Period="1709"; Mora=29509},
    @{Row=47; Period="1708"; Mora=29509},
    @{Row=48; Period="1707"; Mora=29509},
    @{Row=49; Period="1706"; Mora=29509},
    @{Row=50; Period="1705"; Mora=29509},
    @{Row=51; Period="1704"; Mora=29509},
    @{Row=52; Period="1703"; Mora=29509}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.Period   # column E - Periodo Mora
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora      # column F - Valor Mora
}

# The new, wider period/value text reflows the worksheet's "best fit"
# columns (B, C, E, F, G, H, I, J) to slightly larger widths.
$ws.Columns(2).ColumnWidth = 18.54296875
$ws.Columns(3).ColumnWidth = 16.7265625
$ws.Columns(5).ColumnWidth = 13.54296875
$ws.Columns(6).ColumnWidth = 10.1796875
$ws.Columns(7).ColumnWidth = 14.36328125
$ws.Columns(8).ColumnWidth = 19.36328125
$ws.Columns(9).ColumnWidth = 18.08984375
$ws.Columns(10).ColumnWidth = 15
